$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B gets a bit narrower (target stored width ~14.7109375 chars).
# The COM ColumnWidth setter snaps to the nearest 1/6-character pixel grid,
# so 13.86 is the input that lands on the closest reachable width (14.6667).
$ws.Columns.Item(2).ColumnWidth = 13.86

$ws.Range("A1").Value = -0.11327129920958612
$ws.Range("B1").Value = 0.11319305610440011
$ws.Range("A2").Value = -0.091091212324593052
$ws.Range("B2").Value = 0.09085401624787881
$ws.Range("A3").Value = -0.073373377566127829
$ws.Range("B3").Value = 0.073238339803754826
$ws.Range("A4").Value = -0.065238339895820729
$ws.Range("B4").Value = 0.064676900294525552
$ws.Range("A5").Value = -0.061676900341821828
$ws.Range("B5").Value = 0.059751386558663633
$ws.Range("A6").Value = -0.052673566119439741
$ws.Range("B6").Value = 0.052078318837493143
$ws.Range("A7").Value = -0.042078318967049277
$ws.Range("B7").Value = 0.041928461360726921
$ws.Range("A8").Value = -0.03192846149568318
$ws.Range("B8").Value = 0.031646403214322216
$ws.Range("A9").Value = -0.029646403277184596
$ws.Range("B9").Value = 0.029409482569832335
$ws.Range("A10").Value = -0.027409482639368932
$ws.Range("B10").Value = 0.027395026287756252
$ws.Range("A11").Value = -0.024395026368002171
$ws.Range("B11").Value = 0.024368945266354913
$ws.Range("A12").Value = -0.020868945353001767
$ws.Range("B12").Value = 0.020671830009886261
$ws.Range("A13").Value = -0.017171830103025698
$ws.Range("B13").Value = 0.017082425367975773
$ws.Range("A14").Value = -0.0090824255072021742
$ws.Range("B14").Value = 0.0090535825035100714
$ws.Range("A15").Value = -0.0080535825750889245
$ws.Range("B15").Value = 0.0080349973939437547
$ws.Range("A16").Value = -0.0060349974765343539
$ws.Range("B16").Value = 0.0060032244225673992
$ws.Range("A17").Value = -0.0040032245068495342
$ws.Range("B17").Value = 0.0039999998955986271
$ws.Range("A18").Value = -0.016101843840033325
$ws.Range("B18").Value = 0.016090618394127176
$ws.Range("A19").Value = -0.012090618434856371
$ws.Range("B19").Value = 0.012015869701450299
$ws.Range("A20").Value = -0.0080158697453356353
$ws.Range("B20").Value = 0.0080055381094137346
$ws.Range("A21").Value = -0.0040055381537635881
$ws.Range("B21").Value = 0.0039999999552522425
$ws.Range("A22").Value = -0.045704033835097135
$ws.Range("B22").Value = 0.045493490517580781
$ws.Range("A23").Value = -0.040493490579302183
$ws.Range("B23").Value = 0.040097787302288523
$ws.Range("A24").Value = -0.02009778751758784
$ws.Range("B24").Value = 0.019999999781670219
$ws.Range("A25").Value = -0.0050778204878607625
$ws.Range("B25").Value = 0.0050381746332472943
$ws.Range("A26").Value = -0.0025381746870039024
$ws.Range("B26").Value = 0.0024889084455903543
$ws.Range("A27").Value = -0.021373817121840499
$ws.Range("B27").Value = 0.021021155134659519
$ws.Range("A28").Value = -0.01902115518850156
$ws.Range("B28").Value = 0.018799466883339555
$ws.Range("A29").Value = -0.011799466989364404
$ws.Range("B29").Value = 0.0117474010349925
$ws.Range("A30").Value = 0.048252598337438446
$ws.Range("B30").Value = -0.048478541924889651
$ws.Range("A31").Value = -0.014021101091817201
$ws.Range("B31").Value = 0.014000868808237144
$ws.Range("A32").Value = -0.0040008689359982696
$ws.Range("B32").Value = 0.0039999999310307288
